# Fix the value for Bello (column AM) on row 209: it had been stored as the
# text placeholder "####" and needs to be the actual numeric case count.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CasosColombia")

$cell = $ws.Range("AM209")
$cell.Value = 10090

# Move the active selection to the corrected cell (matches the author's
# saved view state).
$cell.Select()
